$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.903.26"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "3.686.98"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "651.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.499"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.146"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.444"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "4.303.60"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "3.662.87"
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("D16").Value = "69.819.09"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "471.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.658"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "3.830.60"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000128"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.167"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").Value = "3.677.50"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "177.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0902"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.932"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000273"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").Value = "  -3.52%  "
$ws.Range("E51").Value = "  -5.23%  "
